{"js": "// The author's commit (\"una cosa peque\u00f1a del doc\" = \"one small thing in the\n// doc\") fixes a small wording mistake in the \"Documentaci\u00f3n\" section:\n// the stray leading \"El \" before \"principalmente\" is removed and the word\n// is capitalized, so \"El principalmente, las llaves primarias...\" becomes\n// \"Principalmente, las llaves primarias...\".\n//\n// (The rest of the underlying XML diff \u2014 namespace-list churn, w:lang\n// locale tags, drawing wp14:anchorId/editId GUIDs, proofErr spell-check\n// markers, bookmark id renumbering, and the TtuloTDC->TtulodeTDC style id\n// rename \u2014 are all re-save artifacts produced when the file was re-saved\n// by a different Word build/platform; they carry no user-visible meaning\n// and are not reproducible through the Word JS API, so we leave them\n// alone and focus on the one real content edit.)\n\nconst searchResults = context.document.body.search(\"El principalmente\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items/text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Expected text \"El principalmente\" was not found in the document.');\n}\n\n// Replace just the matched \"El principalmente\" with \"Principalmente\" \u2014\n// everything that follows (\", las llaves primarias...\") is left untouched.\nsearchResults.items[0].insertText(\"Principalmente\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The author's commit (\"una cosa peque\u00f1a del doc\" = \"one small thing in the\n# doc\") fixes a small wording mistake in the \"Documentaci\u00f3n\" section:\n# the stray leading \"El \" before \"principalmente\" is removed and the word\n# is capitalized, so \"El principalmente, las llaves primarias...\" becomes\n# \"Principalmente, las llaves primarias...\".\n#\n# (The rest of the underlying XML diff -- namespace-list churn, w:lang\n# locale tags, drawing wp14:anchorId/editId GUIDs, proofErr spell-check\n# markers, bookmark id renumbering, and the TtuloTDC->TtulodeTDC style id\n# rename -- are all re-save artifacts produced when the file was re-saved\n# by a different Word build/platform; they carry no user-visible meaning\n# and are not reproducible through the Word object model, so we leave\n# them alone and focus on the one real content edit.)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"El principalmente\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Forward = $true\n$find.Wrap = 0\n\nif ($find.Execute()) {\n    $rng = $find.Parent\n    $rng.Text = \"Principalmente\"\n} else {\n    throw 'Expected text \"El principalmente\" was not found in the document.'\n}\n"}
